$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 6, columns D:F)
$ws.Range("D6").Value = "Action effectuée"
$ws.Range("E6").Value = "Résultat attendu"
$ws.Range("F6").Value = "Résultat obtenu"

# Data rows (rows 7:10, column D only - E and F left blank)
$ws.Range("D7").Value = "Test déplacement pion"
$ws.Range("D8").Value = "Avancer un pion de 1"
$ws.Range("D9").Value = "Avancer un pion de 2 au premier coup"
$ws.Range("D10").Value = "Avancer un pion de 1 puit de 2"

# Style the header row: solid fill (theme 2) + thin border all around
$headerRange = $ws.Range("D6:F6")
$headerRange.Interior.ThemeColor = 2
$headerRange.Interior.Pattern = 1
$headerRange.Borders.LineStyle = 1
$headerRange.Borders.Weight = 2

# Style the body rows: no fill, thin border all around
$bodyRange = $ws.Range("D7:F10")
$bodyRange.Borders.LineStyle = 1
$bodyRange.Borders.Weight = 2

# Autofit the columns used by the table
$ws.Range("D6:F10").Columns.AutoFit()

# Set the active selection to match the target state
$ws.Range("I19").Select()
